# Implementacion de caso de generar recibo de ingresos y egresos
# cuarta iteracion de la implementacion del caso de uso generar recibo
# de ingresos y egresos, se actualiza la plantilla de la iteracion.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# Registrar horas consumidas (columna "Cons.") para la tarea de la fila 7
# (CU Generar recibo de pago.) en el dia 3 -> columna N.
$ws.Range("N7").Value = 1

# Registrar horas consumidas para la tarea de la fila 8
# (CU Generar reporte de ingresos y egresos.) en el dia 5 (columna T)
# y en el dia 6 (columna W).
$ws.Range("T8").Value = 1
$ws.Range("W8").Value = 2

# Actualizar el estatus de la tarea de la fila 11
# (Agregar registro de pago de alumno temporal para el director.)
# de "Hecho" a "En proceso".
$ws.Range("F11").Value = "En proceso"

# Dejar seleccionada la celda correspondiente al ultimo valor capturado.
$ws.Range("W8").Select()
